$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the data in the same order the original author typed it so that the
# shared-strings table (and therefore the <v> indices used by each cell)
# comes out in the same order as the target workbook.
$ws.Range("B1").Value = "Topics "
$ws.Range("C1").Value = "Algorithm"
$ws.Range("D1").Value = "Question "
$ws.Range("B2").Value = "Array"
$ws.Range("C2").Value = "Backtracking, Bit manipulation "
$ws.Range("D2").Value = "Subset II "
$ws.Range("F1").Value = "Solution"
$ws.Range("E1").Value = "Level"
$ws.Range("E2").Value = "Medium"
$ws.Range("A1").Value = "Code"
$ws.Range("G1").Value = "Possible alternatives"
$ws.Range("A2").Value = 90

# Column widths, set to match the auto-fit widths recorded in the target file
# (values chosen so the engine's internal pixel-rounding lands as close as
# possible to the target stored width).
$ws.Columns.Item(3).ColumnWidth = 28.333333333333332
$ws.Columns.Item(4).ColumnWidth = 50.333333333333336
$ws.Columns.Item(6).ColumnWidth = 43.666666666666664
$ws.Columns.Item(7).ColumnWidth = 31.833333333333332

# Page setup / view state.
$ws.PageSetup.Orientation = 1

$ws.Range("E2").Select()
$excel.ActiveWindow.ScrollColumn = 4
